$d = $word.ActiveDocument

# Remove the first 5 paragraphs: "Hi,", "This is with reference...",
# "All the supporting documents...", "Please let me know...", and the
# blank paragraph that follows them — leaving "Regards," as the new
# first paragraph.
for ($i = 0; $i -lt 5; $i++) {
    $p = $d.Paragraphs.Item(1)
    $p.Range.Delete()
}
